# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 432 of the data table
# (pushing the existing rows 432-458 down to 433-459), and populated with
# the new week's data for "Pepino ensalada" at Mercado Mayorista Lo
# Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 432 - this shifts every row at/after 432
# down by one (432 -> 433, ..., 458 -> 459), matching the existing rows'
# data, styles and formatting.
$ws.Rows.Item(432).Insert()

# Populate the newly-inserted row 432 with the new record.
$ws.Cells.Item(432, 1).Value  = 6
$ws.Cells.Item(432, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(432, 3).Value  = "Metropolitana"
$ws.Cells.Item(432, 4).Value  = 44753
$ws.Cells.Item(432, 5).Value  = 13
$ws.Cells.Item(432, 6).Value  = 100112043
$ws.Cells.Item(432, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(432, 8).Value  = "Sin especificar"
$ws.Cells.Item(432, 9).Value  = "Primera"
$ws.Cells.Item(432, 10).Value = 770
$ws.Cells.Item(432, 11).Value = 17000
$ws.Cells.Item(432, 12).Value = 19000
$ws.Cells.Item(432, 13).Value = 18091
$ws.Cells.Item(432, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(432, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(432, 16).Value = 302
$ws.Cells.Item(432, 17).Value = 60
$ws.Cells.Item(432, 18).Value = "Hortaliza"
